# Applies the cryptos.xlsx update described in the commit
# "Updated cryptos list on Wed Nov  8 22:06:53 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '35.776.07'
$ws.Range("E2").Value = '  +0.92%  '
$ws.Range("D3").Value = '1.891.70'
$ws.Range("E3").Value = '  +0.10%  '
$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  -0.31%  '
$ws.Range("B5").Value = 'BNB'
$ws.Range("C5").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D5").Value = '''246.65'
$ws.Range("E5").Value = '  -0.12%  '
$ws.Range("B6").Value = 'XRP'
$ws.Range("C6").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D6").Value = '''0.693'
$ws.Range("E6").Value = '  +0.46%  '
$ws.Range("D7").Value = '''0.999'
$ws.Range("E7").Value = '  -0.20%  '
$ws.Range("D8").Value = '''43.18'
$ws.Range("E8").Value = '  -1.11%  '
$ws.Range("D9").Value = '''56.89'
$ws.Range("E9").Value = '  +9.30%  '
$ws.Range("D10").Value = '''0.357'
$ws.Range("E10").Value = '  +1.24%  '
$ws.Range("D11").Value = '''0.0754'
$ws.Range("E11").Value = '  +1.92%  '
$ws.Range("E12").Value = '  +1.50%  '
$ws.Range("D13").Value = '''14.86'
$ws.Range("E13").Value = '  +12.97%  '
$ws.Range("D14").Value = '''0.790'
$ws.Range("E14").Value = '  +7.42%  '
$ws.Range("D15").Value = '2.162.98'
$ws.Range("E15").Value = '  -0.27%  '
$ws.Range("D16").Value = '''5.06'
$ws.Range("E16").Value = '  +1.92%  '
$ws.Range("D17").Value = '1.879.95'
$ws.Range("E17").Value = '  -0.78%  '
$ws.Range("D18").Value = '35.616.15'
$ws.Range("E18").Value = '  +0.29%  '
$ws.Range("D19").Value = '''73.55'
$ws.Range("E19").Value = '  -0.24%  '
$ws.Range("D20").Value = '0.0₃0831'
$ws.Range("E20").Value = '  +1.00%  '
$ws.Range("D21").Value = '''246.74'
$ws.Range("E21").Value = '  +0.16%  '
$ws.Range("D22").Value = '''13.09'
$ws.Range("E22").Value = '  +1.69%  '
$ws.Range("D23").Value = '''5.19'
$ws.Range("E23").Value = '  +4.65%  '
$ws.Range("E24").Value = '  +3.38%  '
$ws.Range("E25").Value = '  +0.04%  '
$ws.Range("D26").Value = '''2.16'
$ws.Range("E26").Value = '  -1.31%  '
$ws.Range("D27").Value = '''166.29'
$ws.Range("E27").Value = '  -0.03%  '
$ws.Range("D28").Value = '''8.67'
$ws.Range("E28").Value = '  +2.11%  '
$ws.Range("D29").Value = '''18.39'
$ws.Range("E29").Value = '  +0.29%  '
$ws.Range("E30").Value = '  +0.49%  '
$ws.Range("D31").Value = '''4.42'
$ws.Range("E31").Value = '  +4.05%  '
$ws.Range("E32").Value = '  +4.36%  '
$ws.Range("E33").Value = '  +1.45%  '
$ws.Range("E34").Value = '  +17.16%  '
$ws.Range("E36").Value = '  -15.15%  '
$ws.Range("D37").Value = '''0.857'
$ws.Range("E37").Value = '  +0.48%  '
$ws.Range("D38").Value = '''0.0745'
$ws.Range("E38").Value = '  +7.66%  '
$ws.Range("E39").Value = '  -3.40%  '
$ws.Range("D40").Value = '''0.0229'
$ws.Range("E40").Value = '  +6.23%  '
$ws.Range("E41").Value = '  +1.51%  '
$ws.Range("D42").Value = '''16.92'
$ws.Range("E42").Value = '  -1.09%  '
$ws.Range("E43").Value = '  +0.11%  '
$ws.Range("D44").Value = '''14.21'
$ws.Range("E44").Value = '  +17.36%  '
$ws.Range("D45").Value = '1.309.90'
$ws.Range("E45").Value = '  +1.43%  '
$ws.Range("E46").Value = '  -0.87%  '
$ws.Range("D47").Value = '''0.0812'
$ws.Range("E47").Value = '  +0.20%  '
$ws.Range("E48").Value = '  +0.20%  '
$ws.Range("E49").Value = '  -0.32%  '
$ws.Range("E50").Value = '  +0.62%  '
$ws.Range("D51").Value = '''42.64'
$ws.Range("E51").Value = '  -1.66%  '
